$wb = $excel.ActiveWorkbook

function Set-TextValue($range, $text) {
    # Force the value to be stored as literal text (avoids Excel's
    # automatic date / number parsing for strings like "2024-08-17"),
    # while keeping the cell's style/format attribute untouched.
    $range.NumberFormat = "@"
    $range.Value = $text
    $range.ClearFormats()
}

function Update-Sheet($ws) {
    # Bump the "want to go" counter on the first event row.
    $ws.Range("F2").Value = 1017

    # Insert a new row at position 3; this shifts the old rows 3 and 4
    # down to rows 4 and 5 respectively.
    $ws.Rows.Item(3).Insert()

    # The newly inserted row inherits a slightly different auto-created
    # style; restore it to match the other numbering cells (A2, A4, A5)
    # which are bold, centered, top-aligned, with a thin border all
    # around.
    $a3 = $ws.Range("A3")
    $a3.Font.Bold = $true
    $a3.HorizontalAlignment = -4108  # xlCenter
    $a3.VerticalAlignment = -4160    # xlTop
    $a3.Borders.LineStyle = 1

    # Populate the new row 3 with the new event's data.
    $ws.Range("A3").Value = 2
    Set-TextValue $ws.Range("B3") "2024-08-17"
    $ws.Range("C3").Value = "丽水·银泰城次元月稻米同好会（免费入场）"
    $ws.Range("D3").Value = "大洋路大洋河-三区 丽水银泰城"
    $ws.Range("E3").Value = "2024.08.17 15:00-08.17 20:00"
    $ws.Range("F3").Value = 2
    $ws.Range("G3").Value = 25
    $ws.Range("H3").Value = "https://show.bilibili.com/platform/detail.html?id=90624"
    $ws.Range("I3").Value = "//i1.hdslb.com/bfs/openplatform/202408/SFP30Lce1723277030193.jpeg"

    # Renumber the rows that were shifted down by the insert.
    $ws.Range("A4").Value = 3
    $ws.Range("A5").Value = 4
}

for ($i = 1; $i -le $wb.Worksheets.Count; $i++) {
    $ws = $wb.Worksheets.Item($i)
    if ($ws.UsedRange.Rows.Count -ge 4) {
        Update-Sheet $ws
    }
}
